$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value2 = 44277
$ws.Range("J2").Value2 = 150
$ws.Range("K2").Value2 = 11000
$ws.Range("L2").Value2 = 11000
$ws.Range("M2").Value2 = 11000
$ws.Range("P2").Value2 = 550

# Row 4 updates
$ws.Range("D4").Value2 = 44280
$ws.Range("J4").Value2 = 100
$ws.Range("K4").Value2 = 10000
$ws.Range("L4").Value2 = 10000
$ws.Range("M4").Value2 = 10000
$ws.Range("P4").Value2 = 500
